$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.6
$ws.Range("J4").Value = 2.2
$ws.Range("I4").Value = 5.25
$ws.Range("V4").Value = 1.4
$ws.Range("U4").Value = 3
$ws.Range("AD4").Value = 12
$ws.Range("AI4").Value = 17
$ws.Range("AN4").Value = 15

# Row 6
$ws.Range("S6").Value = 3.5
$ws.Range("T6").Value = 1.31
$ws.Range("AS6").Value = 2.14
$ws.Range("AR6").Value = 1.72

# Row 9
$ws.Range("M9").Value = 1.1
$ws.Range("L9").Value = 5
$ws.Range("G9").Value = 1.95
$ws.Range("N9").Value = 7
$ws.Range("K9").Value = 1.91
$ws.Range("I9").Value = 4.33
$ws.Range("R9").Value = 1.5
$ws.Range("Q9").Value = 2.5
$ws.Range("U9").Value = 5
$ws.Range("V9").Value = 1.17
$ws.Range("Y9").Value = 2.2
$ws.Range("Z9").Value = 1.62
$ws.Range("AD9").Value = 17
$ws.Range("AN9").Value = 15
$ws.Range("AS9").Value = 1.9
$ws.Range("AR9").Value = 1.95

# Row 10
$ws.Range("G10").Value = 1.85
$ws.Range("N10").Value = 6.5
$ws.Range("H10").Value = 3
$ws.Range("J10").Value = 2.6
$ws.Range("I10").Value = 5
$ws.Range("K10").Value = 1.95
$ws.Range("M10").Value = 1.11
$ws.Range("R10").Value = 1.5
$ws.Range("T10").Value = 1.24
$ws.Range("Q10").Value = 2.5
$ws.Range("S10").Value = 3.9
$ws.Range("V10").Value = 1.17
$ws.Range("U10").Value = 5
$ws.Range("AE10").Value = 19
$ws.Range("AJ10").Value = 81
$ws.Range("AG10").Value = 6.5
$ws.Range("AC10").Value = 9.5
$ws.Range("AL10").Value = 10
$ws.Range("AF10").Value = 41
$ws.Range("AH10").Value = 6
$ws.Range("AD10").Value = 15
$ws.Range("AS10").Value = 1.98
$ws.Range("AR10").Value = 1.88

# Row 12
$ws.Range("I12").Value = 26
$ws.Range("H12").Value = 7.6
$ws.Range("G12").Value = 1.11
$ws.Range("K12").Value = 2.95
$ws.Range("Z12").Value = 1.42
$ws.Range("AC12").Value = 13
$ws.Range("AI12").Value = 55
$ws.Range("AA12").Value = 6.6
$ws.Range("AF12").Value = 55
$ws.Range("AB12").Value = 5.9
$ws.Range("AE12").Value = 13
$ws.Range("AH12").Value = 19
$ws.Range("Y12").Value = 2.65
$ws.Range("AD12").Value = 6
$ws.Range("AM12").Value = 400
$ws.Range("AN12").Value = 100
$ws.Range("AP12").Value = 700
$ws.Range("AQ12").Value = 400

# Row 13
$ws.Range("L13").Value = 2.27
$ws.Range("J13").Value = 4.9
$ws.Range("I13").Value = 1.7
$ws.Range("H13").Value = 3.8
$ws.Range("G13").Value = 4.65
$ws.Range("R13").Value = 1.9
$ws.Range("O13").Value = 1.28
$ws.Range("P13").Value = 3.45
$ws.Range("AF13").Value = 55
$ws.Range("AL13").Value = 6.7
$ws.Range("AM13").Value = 8.25
$ws.Range("AC13").Value = 16.5
$ws.Range("Z13").Value = 1.88
$ws.Range("AJ13").Value = 90
$ws.Range("AI13").Value = 17.5
$ws.Range("AO13").Value = 13.5
$ws.Range("AA13").Value = 12.5
$ws.Range("AH13").Value = 7.7
$ws.Range("AB13").Value = 29
$ws.Range("AD13").Value = 90

# Row 14
$ws.Range("J14").Value = 5
$ws.Range("L14").Value = 2.25
$ws.Range("I14").Value = 1.65
$ws.Range("AC14").Value = 15

# Row 15
$ws.Range("I15").Value = 3
$ws.Range("G15").Value = 2.35
$ws.Range("AO15").Value = 29

# Row 16
$ws.Range("M16").Value = 1.08
$ws.Range("N16").Value = 8
$ws.Range("AR16").Value = 1.85
$ws.Range("AS16").Value = 2

# Row 17
$ws.Range("G17").Value = 2.15
$ws.Range("L17").Value = 3.25
$ws.Range("J17").Value = 2.88
$ws.Range("I17").Value = 2.8
$ws.Range("T17").Value = 1.85
$ws.Range("S17").Value = 1.95
$ws.Range("AE17").Value = 17
$ws.Range("AC17").Value = 9.5
$ws.Range("AN17").Value = 11
$ws.Range("AO17").Value = 29
$ws.Range("AM17").Value = 17
$ws.Range("AL17").Value = 13

# Row 18
$ws.Range("I18").Value = 2.95
$ws.Range("K18").Value = 2.18
$ws.Range("Q18").Value = 1.62
$ws.Range("R18").Value = 2.02
$ws.Range("V18").Value = 1.42
$ws.Range("H18").Value = 3.4
$ws.Range("L18").Value = 3.45
$ws.Range("Y18").Value = 1.52
$ws.Range("J18").Value = 2.72
$ws.Range("O18").Value = 1.21
$ws.Range("G18").Value = 2.18
$ws.Range("AA18").Value = 9.75
$ws.Range("P18").Value = 3.55
$ws.Range("AB18").Value = 12.5
$ws.Range("Z18").Value = 2.2
$ws.Range("U18").Value = 2.47
$ws.Range("AI18").Value = 11.75
$ws.Range("AQ18").Value = 26
$ws.Range("AO18").Value = 37
$ws.Range("AK18").Value = 250
$ws.Range("AE18").Value = 16
$ws.Range("AL18").Value = 11.5
$ws.Range("AF18").Value = 21
$ws.Range("AG18").Value = 12.5
$ws.Range("AH18").Value = 6.8
$ws.Range("AN18").Value = 10.5
$ws.Range("AP18").Value = 23

# Row 19
$ws.Range("J19").Value = 4.05
$ws.Range("L19").Value = 2.5
$ws.Range("R19").Value = 1.9
$ws.Range("AA19").Value = 11.75
$ws.Range("AF19").Value = 35
$ws.Range("AE19").Value = 32
$ws.Range("AL19").Value = 8
$ws.Range("AK19").Value = 400
$ws.Range("AQ19").Value = 23

# Row 20
$ws.Range("I20").Value = 4.75
$ws.Range("L20").Value = 4.5
$ws.Range("G20").Value = 1.7
$ws.Range("K20").Value = 2.38
$ws.Range("H20").Value = 3.75
$ws.Range("J20").Value = 2.25
$ws.Range("Z20").Value = 2.2
$ws.Range("Y20").Value = 1.62
$ws.Range("AB20").Value = 9.5
$ws.Range("AE20").Value = 13
$ws.Range("AM20").Value = 26

# Row 22
$ws.Range("I22").Value = 3.3
$ws.Range("G22").Value = 2.05
$ws.Range("L22").Value = 3.6
$ws.Range("S22").Value = 1.98
$ws.Range("T22").Value = 1.88
$ws.Range("AO22").Value = 34

# Row 23
$ws.Range("AC23").Value = 12
$ws.Range("AG23").Value = 19
$ws.Range("AQ23").Value = 19

# Row 25
$ws.Range("I25").Value = 1.85
$ws.Range("L25").Value = 2.5
$ws.Range("Q25").Value = 2.05
$ws.Range("R25").Value = 1.75
$ws.Range("AP25").Value = 15

# Row 26
$ws.Range("Q26").Value = 1.73
$ws.Range("R26").Value = 2.08

# Row 27
$ws.Range("N27").Value = 7.7
$ws.Range("P27").Value = 3.5
$ws.Range("V27").Value = 1.38
$ws.Range("X27").Value = 2.8
$ws.Range("U27").Value = 2.82
$ws.Range("AB27").Value = 9.25
$ws.Range("AC27").Value = 8
$ws.Range("AG27").Value = 7.7
$ws.Range("AJ27").Value = 55
$ws.Range("AM27").Value = 23
$ws.Range("AI27").Value = 13.5
$ws.Range("AL27").Value = 12
$ws.Range("AQ27").Value = 37

# Row 28
$ws.Range("O28").Value = 1.25
$ws.Range("U28").Value = 3.25
$ws.Range("R28").Value = 1.95
$ws.Range("V28").Value = 1.33
$ws.Range("P28").Value = 3.75
$ws.Range("Q28").Value = 1.9

# Row 29
$ws.Range("I29").Value = 1.4
$ws.Range("AH29").Value = 9.5
